$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.701.17"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.851.57"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "262.59"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5368"
$ws.Range("E7").Value = "  +2.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3196"
$ws.Range("E8").Value = "  -2.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06982"
$ws.Range("E9").Value = "  +2.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.06"
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7758"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07832"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.846.30"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.66"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.059"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.19"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008022"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.717.31"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.084.65"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.662"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.053"
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.424"
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.223"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.80"
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.700"
$ws.Range("E27").Value = "  +2.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.14"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.54"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.340"
$ws.Range("E30").Value = "  +3.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08759"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.121"
$ws.Range("E32").Value = "  -0.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04886"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7431"
$ws.Range("E34").Value = "  +3.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.145"
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.898"
$ws.Range("E36").Value = "  +1.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.114"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("E38").Value = "  +6.70%  "
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4854"
$ws.Range("E40").Value = "  -0.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9100"
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.48"
$ws.Range("E42").Value = "  -1.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.923"
$ws.Range("E43").Value = "  -2.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.755"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4207"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.125"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1257"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.26"
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05842"
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9016"
$ws.Range("E51").Value = "  +1.11%  "
